# Add the five final cards (rows 43-47) to the bottom of the card data table,
# then move the viewport/selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(41, 10,   0, 15, 0, "yes", "no"),
    @(42,  5, -10,  5, 0, "yes", "no"),
    @(43, 10, -10, 10, 0, "yes", "no"),
    @(44,-10, -10, 10, 0, "yes", "no"),
    @(45,-10,  -5, 10, 0, "no",  "no")
)

$startRow = 43
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}

# Update the view to match the author's final position/selection.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("G48").Select()
